# --------------------------------------------------------------------------
# Speech.docx peer-review pass: the whole body is rewritten (new intro,
# "MotionFit fixes all of this..." paragraph, reshuffled kit/Kinect/targets
# paragraphs, a new closing "Motion Fit is like no other..." + call-to-action
# paragraph, blank spacer paragraphs between sections, proofErr spell/grammar
# markers, and the _GoBack bookmark moving to the end of the new closing
# paragraph).
#
# Rather than a long sequence of Find/Replace calls (error prone for a
# restructuring this size), the new body is supplied as one OOXML fragment
# and dropped in with Range.InsertXML, which is the documented way to swap a
# Range's contents for arbitrary markup in the Word object model.
# --------------------------------------------------------------------------

$d = $word.ActiveDocument

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Fitness and healthy lifestyles are becoming more and more common in todays society, especially with summer coming up, people become more motivated to hit the gym to get that summer body they desperately want. Yet many </w:t></w:r><w:r><w:t>can’t</w:t></w:r><w:r><w:t xml:space="preserve"> afford to get a personal trainer so that just do what they know</w:t></w:r><w:r><w:t xml:space="preserve">. The issue thus arises of </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>gym based</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> injuries. </w:t></w:r><w:r><w:t>The most common cause of injuries, especially regarding gym exercise</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> is lack of proper technique.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>MotionFit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> fixes all of this, this </w:t></w:r><w:r><w:t xml:space="preserve">device </w:t></w:r><w:r><w:t xml:space="preserve">not only provides a vast amount of information about the techniques </w:t></w:r><w:r><w:t>about different</w:t></w:r><w:r><w:t xml:space="preserve"> types of exercise but it also provides a very enhanced feature of tracking your bodily movements as you exercise</w:t></w:r><w:r><w:t>. It</w:t></w:r><w:r><w:t xml:space="preserve"> records your current technique and then shows you how you should be moving to maximize productivity and weight loss and to ensure that you’re exercising correctly to prevent injury.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MotionFit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> kit comes with 10 sensors that you attach to your body when undertaking a particular exercise and that communicate with the device.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Using technologies similar to those found in Nintendo Wii’s and Xbox </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Kinect</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MotionFit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> will provide gym attendees the aid that they need to stay away from injuries in the most effective manner.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">This product targets not only regular gym users, but also those who are maybe avoiding the gym due to lack of knowledge or due to fear of injuries. Furthermore, this device contains state of the art sensors and therefore, even professional athletes can also benefit from using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MotionFit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Motion Fit is like no other, it provides information, education and motivation to any user and can help you reach your dream body without the fear of injury and bring you one step closer today.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/><w:p><w:r><w:t>So get your plans in motion</w:t></w:r><w:r><w:t xml:space="preserve"> and Purchase a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MotionFit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> today to become a fitter, healthier you.</w:t></w:r></w:p><w:p/><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# InsertXML replaces the contents of the range it's called on; $d.Content is
# the whole main story, so this re-writes the entire body in one shot.
# (NOTE: when the very last element before </w:body> is an empty <w:p/>, Word
# folds it into the implicit section-mark paragraph instead of keeping it as
# its own paragraph -- hence one extra trailing <w:p/> in the source markup
# above to end up with the two blank paragraphs the final document needs.)
$d.Content.InsertXML($xml)

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
